$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting
# (values that look numeric would otherwise be auto-converted to numbers)
$ws.Range("D2").Value = "27.104.70"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "1.823.41"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.60"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4632"
$ws.Range("E7").Value = "  -2.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3642"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07291"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8703"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").Value = "1.872.44"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07606"
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.350"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.61"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.477"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008646"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "27.324.26"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.49"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.200"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("D24").Value = "2.086.48"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.67"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.103"
$ws.Range("E28").Value = "  -3.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.27"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.078"
$ws.Range("E30").Value = "  -4.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08920"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.962"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7350"
$ws.Range("E33").Value = "  -3.46%  "
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.140"
$ws.Range("E35").Value = "  -3.31%  "
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.536"
$ws.Range("E37").Value = "  +6.03%  "
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05263"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("E41").Value = "  -2.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.168"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5217"
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1633"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.273"
$ws.Range("E45").Value = "  -3.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4897"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.07"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.636"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06251"
$ws.Range("E51").Value = "  -1.31%  "
